# Fix the "Name: title" line so it carries a closing template delimiter,
# and remove the now-redundant "Broken template key: title" paragraph
# that used to document the unresolved key.

$d = $word.ActiveDocument

# Append the closing "}}" to the template key in the first paragraph.
$d.Content.Find.Execute("Name: title", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Name: title}}", 2)

# Remove the whole "Broken template key: title" paragraph (paragraph 2),
# including its paragraph mark, so the following paragraphs shift up.
$d.Paragraphs.Item(2).Range.Delete()
